$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(9, 2).Value = "0b8815ba976728b9680bc25139d65596"  # 05-050305TC: 0959e04b3aa73c7868f6602bf922730d -> 0b8815ba976728b9680bc25139d65596
$ws.Cells.Item(11, 2).Value = "e7156e7272623db9995be0ecad2c466d"  # 05-050301A: b2c8390815ce162dfbc195a6e3539e5e -> e7156e7272623db9995be0ecad2c466d
$ws.Cells.Item(15, 2).Value = "4c44ff15326426ac613ce87100df7f95"  # 05-050207TP: 23c4bc40f0a8eb34227b73fcade17c44 -> 4c44ff15326426ac613ce87100df7f95
$ws.Cells.Item(17, 2).Value = "cc73dd0015baf47a19781975f7096db1"  # 05-050305TP: 9a0cf9cde071af21b9a8b1250544dbe1 -> cc73dd0015baf47a19781975f7096db1
$ws.Cells.Item(24, 2).Value = "63d3a0ec944dea628d5eca8827e1defa"  # 05-050316TC: 97658fa5e114113e1a449acdf95ddf5c -> 63d3a0ec944dea628d5eca8827e1defa
$ws.Cells.Item(29, 2).Value = "a4b55398550e4d72516019912f9adadd"  # 05-050302A: 880169c4f20521e3d4822a10de244c87 -> a4b55398550e4d72516019912f9adadd
$ws.Cells.Item(34, 2).Value = "0dc8de77eb23c6a495e3e17a2c95a00e"  # 05-050316TP: 3f54be0653ad2244272cdb4c92b66659 -> 0dc8de77eb23c6a495e3e17a2c95a00e
$ws.Cells.Item(121, 2).Value = "7fc21c463874eaab97dd2296ba4cf985"  # 05-050301TP: 1b616b2e73a9c56fefbc2e46caa895da -> 7fc21c463874eaab97dd2296ba4cf985
$ws.Cells.Item(133, 2).Value = "80e9c60bacf8921324725e4400a7339e"  # 05-050312TP: 219aefafdaead1e58e3487a55809ca80 -> 80e9c60bacf8921324725e4400a7339e
$ws.Cells.Item(136, 2).Value = "495222816963cc67b177e8365d512ec9"  # 05-050312TC: 4d8d0cbab3ffe559b044913127f931cb -> 495222816963cc67b177e8365d512ec9
$ws.Cells.Item(159, 2).Value = "307383b7feee0ba305368503ddf9b3b6"  # 05-050203TP: 258d1be4e5ce772f9c17817b83122106 -> 307383b7feee0ba305368503ddf9b3b6
$ws.Cells.Item(162, 2).Value = "7972729b143fdc5bf51ce60b339015ed"  # 05-050308A: 5ade9c4d2c6ee935e6b926f7fb9a0ce9 -> 7972729b143fdc5bf51ce60b339015ed
$ws.Cells.Item(169, 2).Value = "f4cff78a1b3ea628de5dfbe0e7acc5d0"  # 05-050203TC: 934471d5234116c2105632f918393f08 -> f4cff78a1b3ea628de5dfbe0e7acc5d0
$ws.Cells.Item(175, 2).Value = "cfaedcdd888c7e320fa7cd43df31b1d3"  # 05-050303TP: 4d0c7a05dad8d06ddc754c5606b18e82 -> cfaedcdd888c7e320fa7cd43df31b1d3
$ws.Cells.Item(180, 2).Value = "416b9062ba882f4a37aba05fbe3a09d6"  # 05-050303TC: 3b78fbf76c5f265df55a25de18c3e2f9 -> 416b9062ba882f4a37aba05fbe3a09d6
$ws.Cells.Item(183, 2).Value = "ddbb9b1c51ca03aa3190d85516776a80"  # 05-050305A: 3e3a66cbe6076aaf0f431ff00351763d -> ddbb9b1c51ca03aa3190d85516776a80
$ws.Cells.Item(191, 2).Value = "7f4910dc4f551c79753978b33b2f3ee5"  # 05-050314TP: 2660e2641bc201914a9ee3706d1afe79 -> 7f4910dc4f551c79753978b33b2f3ee5
$ws.Cells.Item(198, 2).Value = "7fe3416af2ad8495e7ec4ae8c0caa316"  # 05-050314TC: 64b254efb3909fc569555fa116472ee4 -> 7fe3416af2ad8495e7ec4ae8c0caa316
$ws.Cells.Item(200, 2).Value = "b0a5e21dfda01647e161100001e813a7"  # 05-050306A: 77fc6691c02ede0e98ed5720035b5c68 -> b0a5e21dfda01647e161100001e813a7
$ws.Cells.Item(213, 2).Value = "4c02568b4cef76b34f0f3f3b15a92cdd"  # 05-050303A: 3a717925e8b0b0dcef43f46beb6facf7 -> 4c02568b4cef76b34f0f3f3b15a92cdd
$ws.Cells.Item(227, 2).Value = "440fab7b2c4505575a74c5a8ccba7a93"  # 05-050205TP: 5df12c5655f7fb6f31c94af54215b5aa -> 440fab7b2c4505575a74c5a8ccba7a93
$ws.Cells.Item(228, 2).Value = "7fb108c6761113e297648cd8e77b3d5a"  # 05-050304A: 9fdefb1cd13a71ebba21891c6d2c9ee0 -> 7fb108c6761113e297648cd8e77b3d5a
$ws.Cells.Item(232, 2).Value = "fdad48b520eb8abce493bcc6e1ff2d19"  # 05-050205TC: ee3640aa2c9fca8dbcd22cc7e942fc4a -> fdad48b520eb8abce493bcc6e1ff2d19
$ws.Cells.Item(339, 2).Value = "af2664b8b2fde0d48e5472556ae8eb0f"  # 05-050201TP: 3c91afa877227368cb569ee456c97b0e -> af2664b8b2fde0d48e5472556ae8eb0f
$ws.Cells.Item(420, 2).Value = "bf3569543f5afe0bd329968445d710df"  # 05-0709-070905BTC: 930e9bd628ccd09c643cd2b4a4b8cfad -> bf3569543f5afe0bd329968445d710df
$ws.Cells.Item(464, 2).Value = "abc67da2d08ba146bcdc5fd13e88bc94"  # 05-050204A: cafa73b84464e6ce32c8cccad7acbb7e -> abc67da2d08ba146bcdc5fd13e88bc94
$ws.Cells.Item(465, 2).Value = "222b2740dc4d039f789a3a0fc1ac32f9"  # 05-050313A: 227de680d72f57468721c27f3cc54e37 -> 222b2740dc4d039f789a3a0fc1ac32f9
$ws.Cells.Item(483, 2).Value = "63270cf73239cbc889bcd19902cd5dae"  # 05-050205A: 8e377676ef963f85fc6cdc072adee325 -> 63270cf73239cbc889bcd19902cd5dae
$ws.Cells.Item(485, 2).Value = "ee58cf895ba9ab649fc65b148c27da1f"  # 05-050314A: e300fe9ea0839f8188800edbf88ed7cf -> ee58cf895ba9ab649fc65b148c27da1f
$ws.Cells.Item(506, 2).Value = "c436b06ad587b6de7209d4d37c4d2dfe"  # 05-050202A: 74d987e2cda486e5de1a59d10854a514 -> c436b06ad587b6de7209d4d37c4d2dfe
$ws.Cells.Item(507, 2).Value = "f42ee8538d8c59d7a3f01aea91264041"  # 05-050311A: 444c85f4b5479d65e5f444f1d33ebf48 -> f42ee8538d8c59d7a3f01aea91264041
$ws.Cells.Item(508, 2).Value = "09369acc749d92312a451c120ddfff19"  # 05-050208TP: 4d537e1fa995288b61de8192a7501164 -> 09369acc749d92312a451c120ddfff19
$ws.Cells.Item(513, 2).Value = "3d17facb60c925fe92cf7a83cbe0a0b9"  # 05-050306TP: ad8624bb8862b0276bdeb95a68584b86 -> 3d17facb60c925fe92cf7a83cbe0a0b9
$ws.Cells.Item(521, 2).Value = "8cfa193e3037ecc27a070bebea725fd1"  # 05-050317TC: b53cb95e7b1beed1711de2295117f6fb -> 8cfa193e3037ecc27a070bebea725fd1
$ws.Cells.Item(524, 2).Value = "9611736c43ed545548d0740133a76bf6"  # 05-050203A: 7093e1fa3dcbb0cbb3abfe84b8119398 -> 9611736c43ed545548d0740133a76bf6
$ws.Cells.Item(532, 2).Value = "f3b26f3aca9304a1da95881d5a1f0a05"  # 05-050317TP: a8f9181ed491ed1e0639f790b03e4d96 -> f3b26f3aca9304a1da95881d5a1f0a05
$ws.Cells.Item(555, 2).Value = "f7818e2c4dd111ff6f0f98889b0723cf"  # 05-050201A: 2913280eaeaab28ba119c5ccfd4cc4b2 -> f7818e2c4dd111ff6f0f98889b0723cf
$ws.Cells.Item(580, 2).Value = "519a14e43d540f5d4316a937eceb84d5"  # 05-050308TP: 2e502c7addb80191a57546bebb4ca098 -> 519a14e43d540f5d4316a937eceb84d5
$ws.Cells.Item(624, 2).Value = "38151a0952c7f8a44b10b275fc2f9c73"  # 05-050204TP: 19ad8120ef4e7fd8c61b97404cc3a38f -> 38151a0952c7f8a44b10b275fc2f9c73
$ws.Cells.Item(635, 2).Value = "9d746aeb10115ed1f0f84b0db00bdff4"  # 05-050204TC: 64dc500dba2d19c1084f441cb01c798a -> 9d746aeb10115ed1f0f84b0db00bdff4
$ws.Cells.Item(637, 2).Value = "6634f91fdce01077ee36d8458798247e"  # 05-050302TP: f6a8676f79701259379a58f88f2cf0e1 -> 6634f91fdce01077ee36d8458798247e
$ws.Cells.Item(657, 2).Value = "5d0fc68f08c311936a13669672cd4efc"  # 05-050313TP: ea0bb9282d0b2a34cffce36bf8ed8796 -> 5d0fc68f08c311936a13669672cd4efc
$ws.Cells.Item(663, 2).Value = "05c0ad80709def82a5805cb168d30bd8"  # 05-050313TC: 39ad392d778518bcc663c52f94db70b2 -> 05c0ad80709def82a5805cb168d30bd8
$ws.Cells.Item(673, 2).Value = "37e8eefe5053e680a6759078e74d2ad2"  # 05-050208A: cbb5f3ebf4381d6e4b27c30867ccb7f7 -> 37e8eefe5053e680a6759078e74d2ad2
$ws.Cells.Item(674, 2).Value = "d1c2d7f30357d105c5d61fb44373b1e1"  # 05-050317A: ebca48fdbfb7ccaf67e04147f6865b4e -> d1c2d7f30357d105c5d61fb44373b1e1
$ws.Cells.Item(688, 2).Value = "aa8738872bea2ef2a790108b17af6b05"  # 05-050206TP: 15158a0991e3dad4fd94dfa5f9c8f3aa -> aa8738872bea2ef2a790108b17af6b05
$ws.Cells.Item(693, 2).Value = "1a4c709560f9ac74cd16d4d06f649f67"  # 05-050206TC: ebe45a973afff04c51d23b1b99035c84 -> 1a4c709560f9ac74cd16d4d06f649f67
$ws.Cells.Item(708, 2).Value = "6da206dbbe2cd535071e418874d32ea6"  # 05-050304TC: 1f4e61800299458a2b76285fe27abd7a -> 6da206dbbe2cd535071e418874d32ea6
$ws.Cells.Item(711, 2).Value = "96d3fef023b15c9f53721cbcb8462f8b"  # 05-050206A: 04461bccc6ab0a10df8f8af8fdc52745 -> 96d3fef023b15c9f53721cbcb8462f8b
$ws.Cells.Item(712, 2).Value = "5253c5919fcaa57059b6fd0e62b6e482"  # 05-050315A: 9866185052e14f49b301a47e90057f55 -> 5253c5919fcaa57059b6fd0e62b6e482
$ws.Cells.Item(723, 2).Value = "98019ef0f7fb12853f6b8ca47fcfa2fc"  # 05-050304TP: b45340bd18cd2b4943af8829769651fb -> 98019ef0f7fb12853f6b8ca47fcfa2fc
$ws.Cells.Item(737, 2).Value = "8dab9df0ec688cdbc74d3e8a5402ff29"  # 05-050316A: 49281e820c63918dbaceddd9728ab270 -> 8dab9df0ec688cdbc74d3e8a5402ff29
$ws.Cells.Item(741, 2).Value = "43c26e689348d34beec01d70727391ef"  # 05-050207A: fd03ec2e714e596c6312367eb6d1c042 -> 43c26e689348d34beec01d70727391ef
$ws.Cells.Item(750, 2).Value = "d0b2277be333de2452ed30e534a8a8b0"  # 05-050315TP: 4c1553eee3fd1eb9927e78dac8b8963e -> d0b2277be333de2452ed30e534a8a8b0
$ws.Cells.Item(827, 2).Value = "296def4cb8a958d7f22c305912bbdeb2"  # 05-050202TP: c535bd182261cc93be3c4531f608bc46 -> 296def4cb8a958d7f22c305912bbdeb2
$ws.Cells.Item(838, 2).Value = "4df54f854c4d53c1fde24fe757645962"  # 05-050311TC: 71f8b444f7700ac0320c268e6589b6c9 -> 4df54f854c4d53c1fde24fe757645962
$ws.Cells.Item(843, 2).Value = "b89027a066248635e7690c9717ed06ec"  # 05-050311TP: 3d731832fb79f3cbf265acdce71ca60f -> b89027a066248635e7690c9717ed06ec
$ws.Cells.Item(862, 2).Value = "71b3057b77a4c1d60dc4f6e0309a33c7"  # 05-050309TC: 8c360e20f2851665840633e15dbd912e -> 71b3057b77a4c1d60dc4f6e0309a33c7
$ws.Cells.Item(882, 2).Value = "d878f735a89572d2273c1e98708e28dd"  # 03-030032A: c9c849f03081bb7a17b5eba5feebb7ea -> d878f735a89572d2273c1e98708e28dd
